$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 (Klinikum Merzig gGmbH entry) - update hospital coordinates, travel time and distance
$ws.Range("E9").Value = 49.4590255
$ws.Range("F9").Value = 6.6299114
$ws.Range("G9").Value = 12.66
$ws.Range("H9").Value = 11.26

# Row 10 (Klinikum Merzig gGmbH entry) - update hospital coordinates, travel time and distance
$ws.Range("E10").Value = 49.4590255
$ws.Range("F10").Value = 6.6299114
$ws.Range("G10").Value = 25.12
$ws.Range("H10").Value = 16.5

$wb.Save()
